$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "The applicant has experience working on projects that involve ReactJS, JavaScript, Python, and some aspects of web development, which align well with the company's requirements. The applicant has also worked on projects related to data analysis, computer vision, and image processing, showcasing a strong technical background. However, the applicant lacks direct experience with MongoDB and NodeJS, which are specifically mentioned in the job description. MongoDB is a crucial skill required for the role, and the applicant's lack of experience in this area may have impacted the score. While the applicant possesses valuable technical skills, the absence of direct experience with the required technologies seems to have influenced the score."
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = "The applicant demonstrates a strong willingness to work in Japan, expressing an interest in the country's work culture and the opportunity to learn a new language. The applicant also acknowledges a need for time to adapt to a new culture and demonstrates a desire to contribute to cutting-edge AI/ML technologies during the internship. Additionally, the applicant's acknowledgment of strengths in teamwork and problem-solving, along with a proactive approach to addressing weaknesses and past challenges, reflects qualities of adaptability and teamwork. Therefore, based on the provided answers and the alignment with the company's required soft skills, the applicant is scored 4 out of 5 for personality and willingness to work in Japan."

$ws.Range("C3").Value = "The applicant has a good score because the projects they have worked on align well with the job description. The 'SHAMIYANA APP' project showcases their web development skills using Dart and Firebase, while the 'RAPID' project highlights their expertise in NodeJS, ExpressJS, Socket.IO, and WebRTC, which are in line with the required skills for the job. Additionally, the 'SMART SENSING MIDDLEWARE' project demonstrates their proficiency in ReactJS, further matching the company's requirements. The applicant possesses a strong skill set that covers most of the required technologies and demonstrates practical experience in web and mobile application development, making them a suitable candidate for the SDE Intern role."
$ws.Range("E3").Value = "The applicant demonstrates a strong willingness to work in Japan, evident from their deep interest in the field of AI and their admiration for Japan's work culture. Their openness to adapt to new cultures, as well as their motivation to learn a new language while working in Japan, reflects their flexibility and adaptability. The applicant's emphasis on the importance of teamwork and their acknowledgment of areas for improvement, along with their ability to overcome challenges through collaborative problem-solving, indicate good communication and teamwork skills. Overall, the applicant's responses highlight a positive attitude, adaptability, and a strong interest in being a part of Japan's work culture, warranting a score of 4."

$ws.Range("C4").Value = "The applicant has been given a score of 74.67 because the candidate has demonstrated strong technical skills in projects such as 'Cloudphysician's Vital Extraction Challenge' involving Computer Vision and Image Processing, 'Multi Model Data Analysis for Annotation of Human Activities' with the use of Python, PyTorch, Tensorflow, and Sklearn, and 'LLMGuard' which involved machine learning models trained for detecting hate speech and bias. Although the applicant lacks direct experience with MERN stack technologies, the demonstrated skills in Python, PyTorch, and ReactJS showcase a strong foundation in programming and machine learning, indicating the potential for quick learning and adaptation to the technologies required by the company. The diverse project experience also indicates the ability to work on complex projects and analytical problem-solving skills."
$ws.Range("E4").Value = "The applicant has expressed a keen interest in working in Japan, stating admiration for the country's work culture, a willingness to learn a new language, and an openness to adapting to new cultures. The applicant's responses also demonstrate good communication skills, an inclination towards teamwork, and problem-solving abilities, as evident from the example of overcoming project setbacks and learning from the experience. The applicant also expresses a desire for personal and professional growth during the internship, indicating an adaptable and focused attitude. However, the hesitancy to work alone and the need to improve in acquiring new skills may require some support and mentorship in a new work environment. Overall, the applicant shows strong potential for fitting into the company's culture and working effectively in a cross-cultural setting."
